$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text
# (matching the original inlineStr cell type in the sheet).
$textForcedRefs = @(
    "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D32", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D45", "D47", "D49", "D50", "D51"
)

foreach ($ref in $textForcedRefs) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
}

# Apply the updated price / volume values from the latest crypto data pull.
$ws.Range("D2").Value = "60.993.79"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.883.06"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "587.70"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "140.26"
$ws.Range("E6").Value = "  -4.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.491"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "6.88"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("D11").Value = "0.429"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").Value = "0.0000217"
$ws.Range("E12").Value = "  -3.96%  "
$ws.Range("D13").Value = "32.31"
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "3.365.16"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "61.058.03"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "2.888.10"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "6.51"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").Value = "425.83"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "13.27"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "0.652"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "6.92"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "79.87"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "10.48"
$ws.Range("E24").Value = "  -4.74%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "2.05"
$ws.Range("E26").Value = "  -6.50%  "
$ws.Range("D27").Value = "11.34"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  -8.46%  "
$ws.Range("D30").Value = "6.71"
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "25.82"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").Value = "0.0₃0853"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "0.971"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("D36").Value = "5.44"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").Value = "49.07"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Value = "2.80"
$ws.Range("E38").Value = "  -7.47%  "
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "8.32"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").Value = "40.16"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "0.264"
$ws.Range("E43").Value = "  -7.32%  "
$ws.Range("D44").Value = "2.656.08"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "132.71"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "343.21"
$ws.Range("E47").Value = "  -9.80%  "
$ws.Range("D49").Value = "22.58"
$ws.Range("E49").Value = "  -5.32%  "
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "1.93"
$ws.Range("E51").Value = "  -3.70%  "

# Restore the original (default) style on cells where we forced text format,
# so only the cell value changes and formatting stays as it was.
foreach ($ref in $textForcedRefs) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Write-Host "Applied crypto price/volume updates"
